# Apply the committed changes to the GAIL (INDIA) Ltd equity history sheet:
#  1. Change the number format of the existing "trade_date" column (I2:I357)
#     from date-only (YYYY-MM-DD) to date-time (YYYY-MM-DD HH:MM:SS).
#  2. Append 4 new daily rows (358-361) of trade data, whose "trade_date"
#     column (I) keeps the original date-only format, matching the style the
#     rest of the column used to have before step 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-format the existing trade_date column (rows 2-357) ---
$dateTimeRange = $ws.Range("I2:I357")
$dateTimeRange.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- 2. Append the new rows of data ---
$newRows = @(
    @{ Row = 358; A = 159.98;  B = 46049; D = 161.4;  E = 158.37; F = 160.85; H = 20747970; I = 46049 },
    @{ Row = 359; A = 168.14;  B = 46050; D = 168.65; E = 161.24; F = 161.35; H = 11179417; I = 46050 },
    @{ Row = 360; A = 167.38;  B = 46051; D = 171.87; E = 165.76; F = 168;    H = 15182093; I = 46051 },
    @{ Row = 361; A = 167.29;  B = 46052; D = 168.05; E = 165.25; F = 166.38; H = 15828320; I = 46052 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 3).Value = "NSE"
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = "GAIL"
    $ws.Cells.Item($row, 8).Value = $r.H

    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 9).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 10).Value = "INE129A01019"
    $ws.Cells.Item($row, 11).Value = "GAIL (INDIA) Ltd"
    $ws.Cells.Item($row, 12).Value = "GAIL"
    $ws.Cells.Item($row, 13).Value = "BREEZE"
}
